$wb = $excel.ActiveWorkbook

# --- Rushing sheet (Week 16 actuals logged for J.Jackson, J.Kelley, L.Rountree) ---
$rushing = $wb.Worksheets.Item("Rushing")

# J.Jackson (row 3)
$rushing.Range("C3").Value = 105
$rushing.Range("D3").Value = 57
$rushing.Range("E3").Value = 10
$rushing.Range("F3").Value = 40

# J.Kelley (row 4)
$rushing.Range("C4").Value = 22
$rushing.Range("D4").Value = 18
$rushing.Range("E4").Value = 3
$rushing.Range("F4").Value = 14

# L.Rountree (row 5)
$rushing.Range("C5").Value = 15
$rushing.Range("D5").Value = 7
$rushing.Range("E5").Value = 6
$rushing.Range("F5").Value = 5

# --- Receiving sheet (Week 16 actuals logged for J.Jackson, J.Kelley, L.Rountree) ---
$receiving = $wb.Worksheets.Item("Receiving")

# J.Jackson (row 2)
$receiving.Range("C2").Value = 66
$receiving.Range("D2").Value = 53
$receiving.Range("E2").Value = 4
$receiving.Range("F2").Value = 2
$receiving.Range("G2").Value = 14
$receiving.Range("H2").Value = 12

# J.Kelley (row 3)
$receiving.Range("C3").Value = 9
$receiving.Range("D3").Value = 8
$receiving.Range("E3").Value = 0
$receiving.Range("F3").Value = 0

# L.Rountree (row 4)
$receiving.Range("C4").Value = 4
$receiving.Range("D4").Value = 3
$receiving.Range("E4").Value = 0
$receiving.Range("F4").Value = 0
